$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.608.62'
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").Value = '2.373.18'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D5").Value = "'0.674"
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").Value = "'239.44"
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("D7").Value = "'73.60"
$ws.Range("E7").Value = '  +6.42%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +19.45%  '
$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = '  +7.45%  '
$ws.Range("D11").Value = "'29.84"
$ws.Range("E11").Value = '  +11.70%  '
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").Value = '2.723.07'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = "'16.95"
$ws.Range("E14").Value = '  +7.17%  '
$ws.Range("D15").Value = "'6.78"
$ws.Range("E15").Value = '  +8.45%  '
$ws.Range("D16").Value = "'0.900"
$ws.Range("E16").Value = '  +6.13%  '
$ws.Range("D17").Value = '2.371.43'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '44.594.04'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("E19").Value = '  +5.08%  '
$ws.Range("D20").Value = "'77.60"
$ws.Range("E20").Value = '  +4.62%  '
$ws.Range("D21").Value = "'6.49"
$ws.Range("E21").Value = '  +3.63%  '
$ws.Range("D22").Value = "'255.45"
$ws.Range("E22").Value = '  +2.62%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -3.82%  '
$ws.Range("D25").Value = "'2.53"
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").Value = "'10.45"
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("D27").Value = "'2.30"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").Value = "'22.54"
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = "'1.60"
$ws.Range("E29").Value = '  +5.52%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'174.28"
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("E32").Value = '  +5.63%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0743"
$ws.Range("E33").Value = '  +7.43%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'5.22"
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = '  +3.94%  '
$ws.Range("D36").Value = "'3.94"
$ws.Range("E36").Value = '  +8.11%  '
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").Value = "'6.55"
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = "'0.0273"
$ws.Range("E39").Value = '  +6.73%  '
$ws.Range("D40").Value = "'20.28"
$ws.Range("E40").Value = '  +10.62%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = "'8.88"
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("E43").Value = '  +3.60%  '
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").Value = "'0.0984"
$ws.Range("E45").Value = '  +3.35%  '
$ws.Range("D46").Value = "'98.86"
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = "'0.184"
$ws.Range("E47").Value = '  +12.23%  '
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = "'4.47"
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("D50").Value = '1.445.52'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '2.597.64'
$ws.Range("E51").Value = '  -0.07%  '
